$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "party_pwib" column description row at the bottom of the table
$ws.Range("A79").Value = "party_pwib"
$ws.Range("B79").Value = "Partyname in PWIB policy position data"

# Reflect the updated selection/active cell from the authored workbook
$ws.Range("B80").Select()
